# Update "想去人数" (column F) values across the worksheets to reflect the
# latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet order per workbook.xml: 1=展览, 2=演出, 3=本地生活, 4=全部类型
$wsExhibition = $wb.Worksheets.Item(1)
$wsShow       = $wb.Worksheets.Item(2)
$wsAllTypes   = $wb.Worksheets.Item(4)

# --- 展览 (Exhibition) sheet: column F updates by row ---
$exhibitionUpdates = @{
    2  = 4883
    3  = 2763
    5  = 2821
    7  = 1956
    9  = 1736
    10 = 755
    11 = 485
    12 = 243
    13 = 420
    14 = 1070
    15 = 302
    16 = 7
    17 = 92
    18 = 74
    19 = 703
    20 = 703
    23 = 665
    24 = 750
    26 = 6
    27 = 42
    28 = 532
    29 = 39
    31 = 1594
    32 = 373
    33 = 44
    34 = 1523
    35 = 204
    36 = 2363
    37 = 400
    39 = 618
    41 = 65
    43 = 799
    44 = 1499
    45 = 218
    47 = 497
    48 = 51
    49 = 80
    50 = 113
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# --- 演出 (Show) sheet: column F updates by row ---
$showUpdates = @{
    4  = 95
    12 = 42
}

foreach ($row in $showUpdates.Keys) {
    $wsShow.Cells.Item($row, 6).Value = $showUpdates[$row]
}

# --- 全部类型 (All Types) sheet: column F updates by row ---
$allTypesUpdates = @{
    2  = 4883
    3  = 2763
    4  = 2821
    5  = 1736
    7  = 755
    8  = 485
    9  = 243
    10 = 420
    11 = 1070
    12 = 302
    13 = 92
    14 = 74
    15 = 704
    16 = 704
    18 = 665
    19 = 750
    21 = 95
    22 = 95
    25 = 42
    26 = 532
    27 = 39
    29 = 1594
    30 = 373
    31 = 44
    33 = 2363
    34 = 400
    38 = 42
    39 = 618
    41 = 65
    43 = 799
    44 = 1499
    46 = 218
    47 = 497
    48 = 51
    49 = 80
    50 = 113
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
